$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column T ("LGD_value") before the existing Amortization_type column ---
$ws.Columns("T:T").Insert()

# Header
$ws.Range("T1").Value = "LGD_value"

# Data rows 2-5 (fraction values)
$ws.Range("T2").Value = 0.35
$ws.Range("T3").Value = 0.4
$ws.Range("T4").Value = 0.4
$ws.Range("T5").Value = 0.4

# Fix the renamed shared string that landed in U5 ("I-FINE" -> "I_FINE")
$ws.Range("U5").Value = "I_FINE"

# --- Cell format tweak: F2 becomes left-aligned ---
$ws.Range("F2").HorizontalAlignment = -4131

# --- Add new row 6 of data ---
$ws.Rows("6:6").Insert()

$ws.Range("A6").Value = "ABC200"
$ws.Range("B6").Value = "5-"
$ws.Range("C6").Value = "5-"
$ws.Range("D6").Value = "5-"
$ws.Range("E6").Value = "5-"
$ws.Range("F6").Value = "5-"

$ws.Range("G5:H5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = "100000002202024"
$ws.Range("H6").Value = "789000002202024 X"

$ws.Range("I6").Value = "ABC200"
$ws.Range("J6").Value = "UNSECURED"
$ws.Range("K6").Value = "01 - Stage 1"

$ws.Range("M5:N5").Copy()
$ws.Range("M6").PasteSpecial(-4122)
$ws.Range("M6").Value = 45808
$ws.Range("N6").Value = 45810

$ws.Range("P6").Value = "A3B"
$ws.Range("Q6").Value = "A3B"

$ws.Range("T6").Value = 0.45
$ws.Range("U6").Value = "I_FINE"
$ws.Range("V6").Value = 138
$ws.Range("W6").Value = "B"

$ws.Range("X5").Copy()
$ws.Range("X6").PasteSpecial(-4122)
$ws.Range("X6").Value = 45808

$ws.Range("Y6").Value = "EUR"
$ws.Range("Z6").Value = 176
$ws.Range("AA6").Value = 0.4

# --- Column width touch-ups (best-effort; engine rounds to whole characters) ---
$ws.Columns("G:G").ColumnWidth = 17
$ws.Columns("H:H").ColumnWidth = 18
$ws.Columns("K:K").ColumnWidth = 10
$ws.Columns("M:N").ColumnWidth = 10
$ws.Columns("X:X").ColumnWidth = 10

# --- View settings: zoom 85%, reset scroll so topLeftCell goes back to A1, select G2 ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.Zoom = 85
$ws.Range("G2").Select()
